$wb = $excel.ActiveWorkbook

$tenMst = $wb.Worksheets.Item("TenMst")

# 1) Insert the new hidden "SystemConf" sheet right after "TenMst" so that
#    the new shared strings it introduces (grp_cd, grp_eda_no, val, param,
#    biko) land before "ItemCd3" (added later on TenMst), matching the
#    order new strings were appended to xl/sharedStrings.xml.
$systemConf = $wb.Worksheets.Add($null, $tenMst)
$systemConf.Name = "SystemConf"

$systemConf.Range("A1").Value = "hp_id"
$systemConf.Range("B1").Value = "grp_cd"
$systemConf.Range("C1").Value = "grp_eda_no"
$systemConf.Range("D1").Value = "val"
$systemConf.Range("E1").Value = "param"
$systemConf.Range("F1").Value = "biko"
$systemConf.Range("G1").Value = "create_date"
$systemConf.Range("H1").Value = "create_id"
$systemConf.Range("I1").Value = "create_machine"
$systemConf.Range("J1").Value = "update_date"
$systemConf.Range("K1").Value = "update_id"
$systemConf.Range("L1").Value = "update_machine"

$systemConf.Range("A2").Value = 22
$systemConf.Range("B2").Value = 2008
$systemConf.Range("C2").Value = 0
$systemConf.Range("D2").Value = 0
$systemConf.Range("G2").Value = 44450.687767337964
$systemConf.Range("G2").NumberFormat = "mm:ss.0"
$systemConf.Range("H2").Value = 0
$systemConf.Range("J2").Value = 44450.687767337964
$systemConf.Range("J2").NumberFormat = "mm:ss.0"
$systemConf.Range("K2").Value = 0

$systemConf.Range("A2").Select() | Out-Null
$systemConf.Visible = [Microsoft.Office.Interop.Excel.XlSheetVisibility]::xlSheetHidden

# 2) Append a new row 5 to "TenMst", cloned from row 3 with a few cells
#    changed (item code, kohatu_kbn and yj_cd).
$tenMst.Range("A3:GG3").Copy()
$tenMst.Range("A5:GG5").PasteSpecial()
$tenMst.Range("B5").Value = "ItemCd3"
$tenMst.Range("DB5").Value = 4
$tenMst.Range("DX5").Value = "test"
$tenMst.Range("DE11").Select() | Out-Null

# 3) Re-fetch the remaining sheets now that the collection has been
#    changed (indices shift once "SystemConf" is inserted) and update
#    the view state: "YohoSetMst" becomes the active/selected sheet
#    (it was "IpnKasanExcludeItem" before).
$yohoSetMst = $wb.Worksheets.Item("YohoSetMst")
$ipnKasanExcludeItem = $wb.Worksheets.Item("IpnKasanExcludeItem")

$ipnKasanExcludeItem.Range("B5").Select() | Out-Null

$yohoSetMst.Range("A2").Select() | Out-Null
$yohoSetMst.Activate() | Out-Null
